# Update "horarios" (schedule) workbook with the latest scrape results.
# Última actualización (last update) timestamp moves from 01:50:13 to 02:15:29,
# and the upcoming-arrivals table is refreshed accordingly on each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912" ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:15:29"

$ws1.Range("A6").Value = "02:15:29"
$ws1.Range("B6").Value = "02:58"
$ws1.Range("C6").Value = "215_ALUAR"
$ws1.Range("D6").Value = 43
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = "02:15:29"
$ws1.Range("B7").Value = "03:56"
$ws1.Range("C7").Value = "14_ABASTO"
$ws1.Range("D7").Value = 101
$ws1.Range("E7").Value = "LP1912"

$ws1.Range("A8").Value = "02:15:29"
$ws1.Range("B8").Value = "04:01"
$ws1.Range("C8").Value = "81_EL PELIGRO"
$ws1.Range("D8").Value = 106
$ws1.Range("E8").Value = "LP1912"

# --- Sheet "LP1912-215" ---
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:15:29"

$ws2.Range("A6").Value = "02:15:29"
$ws2.Range("B6").Value = "02:58"
$ws2.Range("C6").Value = "215_ALUAR"
$ws2.Range("D6").Value = 43
$ws2.Range("E6").Value = "LP1912"

# --- Sheet "6203-6173" ---
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 02:15:29"
